# Rename the first person from "Terrence1" to "Terrence2" and point their
# hyperlink at the new e-mail address, then drop the other three rows'
# data (keeping the "Hyperlink" cell style in column B) and extend the
# sheet with a few new blank styled rows further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name in A2.
$ws.Range("A2").Value2 = "Terrence2"

# Remove every existing hyperlink on the sheet (there's no per-item delete
# in this object model, only a blanket one scoped to the whole sheet).
$ws.Range("B2").Hyperlinks.Delete()

# Re-create just the one hyperlink that should remain, in place on B2, and
# keep its original "Hyperlink" cell style untouched by mutating the
# (now-empty) hyperlink reference instead of calling Hyperlinks.Add().
$ws.Range("B2").Hyperlinks.Item(1).Address = "mailto:terrencereinhardt2@gmail.com"
$ws.Range("B2").Value2 = "terrencereinhardt2@gmail.com"

# Clear out the old data for rows 3-5 (names in column A disappear
# entirely; column B keeps its "Hyperlink" style but becomes empty).
$ws.Range("A3:A5").ClearContents()
$ws.Range("B3:B5").ClearContents()

# Add three new blank rows (7, 8, 9) in column B that carry the same
# "Hyperlink" style as the rest of column B, extending the used range.
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"

# Update the active selection to B2.
$ws.Range("B2").Select() | Out-Null
